# Update input_parameters.xlsx sheet ("update after Manas check")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 5 (Vegetation / treatment / KBDI dummy-data column) ----
$ws.Range("E5").Value = "none"
$ws.Range("F5").Value = "FRB=Fuel reduction burn"

# ---- Row 6 (KBDI) ----
$ws.Range("D6").Value = "0-200"
$ws.Range("E6").Value = 150

# ---- Row 7 (w) ----
$ws.Range("D7").ClearContents()
$ws.Range("E7").Value = 1

# ---- Row 8 (Elevation / 3D-Relief) ----
$ws.Range("E8").Value = "flat"

# ---- Row 18 (Properties bounds) ----
$ws.Range("C18").Value = "Bounds (North, East, South, West)"

# ---- New explanatory notes below the table (rows 23-33) ----
$ws.Range("B23").Value = "For the parameter w use ""1""."

$ws.Range("B24").Value = "FRB means if there was fuel reduction burn (FRB) applied to the landscape or not. Because if yes, it reduced the spread of fire. For now use ""0"" meaning no FRB was applied anywhere."
$ws.Range("B24").Font.Name = "Courier New"
$ws.Range("B24").Font.Size = 8
$ws.Range("B24").Font.Color = 0

$ws.Range("B25").Value = "For the KBDI use 150 (see details below)"
$ws.Range("B25").Font.Name = "Courier New"
$ws.Range("B25").Font.Size = 8
$ws.Range("B25").Font.Color = 0

$ws.Range("B27").Value = "KBDI = 0 - 50: Soil moisture and large class fuel moistures are high and do not contribute much to fire intensity. Typical of spring dormant season following winter precipitation."
$ws.Range("B27").Font.Name = "Courier New"
$ws.Range("B27").Font.Size = 8
$ws.Range("B27").Font.Color = 0

$ws.Range("B29").Value = "KBDI = 50 - 100: Typical of late spring, early growing season. Lower litter and duff layers are drying and beginning to contribute to fire intensity."
$ws.Range("B29").Font.Name = "Courier New"
$ws.Range("B29").Font.Size = 8
$ws.Range("B29").Font.Color = 0

$ws.Range("B31").Value = "KBDI = 100 - 150: Typical of late summer, early fall. Lower litter and duff layers actively contribute to fire intensity and will burn actively."
$ws.Range("B31").Font.Name = "Courier New"
$ws.Range("B31").Font.Size = 8
$ws.Range("B31").Font.Color = 0

$ws.Range("B33").Value = "KBDI = 150 - 200: Often associated with more severe drought with increased wildfire occurrence. Intense, deep burning fires with significant downwind spotting can be expected. Live fuels can also be expected to burn actively at these levels."
$ws.Range("B33").Font.Name = "Courier New"
$ws.Range("B33").Font.Size = 8
$ws.Range("B33").Font.Color = 0

# ---- View state: zoom to 76%, select E13 ----
[void]$ws.Range("E13").Select()
$excel.ActiveWindow.Zoom = 76
